# Insert a new data row at sheet row 457 (pushing existing rows 457-510 down to 458-511)
# and populate the newly inserted row with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 457; Excel shifts rows 457..510 down to 458..511
# and extends formatting (e.g. the date style on column D) into the new row.
$ws.Rows.Item(457).Insert()

# Populate the new row 457 with the new record values.
$ws.Range("A457").Value = 10
$ws.Range("B457").Value = "Vega Modelo de Temuco"
$ws.Range("C457").Value = "La Araucanía"
$ws.Range("D457").Value = 45124
$ws.Range("E457").Value = 9
$ws.Range("F457").Value = 100112017
$ws.Range("G457").Value = "Apio"
$ws.Range("H457").Value = "Americana (o)"
$ws.Range("I457").Value = "Primera"
$ws.Range("J457").Value = 250
$ws.Range("K457").Value = 8000
$ws.Range("L457").Value = 10000
$ws.Range("M457").Value = 8400
$ws.Range("N457").Value = "$/docena de matas"
$ws.Range("O457").Value = "Provincia del Elquí"
$ws.Range("P457").Value = 1400
$ws.Range("Q457").Value = 6
$ws.Range("R457").Value = "Hortaliza"
